# Auto refresh - 16-02-2026 11:39:09.39
#
# The "Excel_vs_ML" sheet's meta/diagnostic columns (P..W) are reshuffled:
#   - a new "DSP_meta" column is introduced (copied from column B)
#   - Total_Budget_meta / Flight_Start_Date_meta / Flight_End_Date_meta are
#     reordered so the two date columns come first
#   - the old "Campaign_Status" column is dropped and ML_Prediction moves
#     into its place
#   - the old "Excel_vs_ML_Disagree" column (last, W) is removed entirely
#
# The "Exec_Summary" sheet loses its "Excel vs ML Disagreement Count" row
# and the LAST_REFRESH_UTC timestamp is bumped.

$wb = $excel.ActiveWorkbook
$dateFmt = "YYYY-MM-DD HH:MM:SS"

# ---------------------------------------------------------------------
# Sheet 1: Excel_vs_ML
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Excel_vs_ML")

# Find the last used row on the sheet (header is row 1, data starts row 2).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $oldTotalBudgetMeta = $ws.Range("P$r").Value2      # Total_Budget_meta
    $oldFlightStartMeta = $ws.Range("Q$r").Value2      # Flight_Start_Date_meta
    $oldFlightEndMeta   = $ws.Range("R$r").Value2      # Flight_End_Date_meta
    $oldMlPrediction    = $ws.Range("S$r").Value2      # ML_Prediction
    $dsp                = $ws.Range("B$r").Value2      # DSP

    # New P: Flight_Start_Date_meta (date-formatted)
    $ws.Range("P$r").NumberFormat = $dateFmt
    $ws.Range("P$r").Value2 = $oldFlightStartMeta

    # New Q: Flight_End_Date_meta (date-formatted)
    $ws.Range("Q$r").NumberFormat = $dateFmt
    $ws.Range("Q$r").Value2 = $oldFlightEndMeta

    # New R: Total_Budget_meta (plain number, no special format)
    $ws.Range("R$r").ClearFormats()
    $ws.Range("R$r").Value2 = $oldTotalBudgetMeta

    # New S: DSP_meta (copied from column B)
    $ws.Range("S$r").Value2 = $dsp

    # New T: ML_Prediction (was column S; replaces old Campaign_Status)
    $ws.Range("T$r").Value2 = $oldMlPrediction

    # U (Budget_At_Risk) and V (ML_Early_Warning) stay put untouched.
}

# Header row: same reshuffle as the data rows, plus the new column name.
$ws.Range("P1").Value2 = "Flight_Start_Date_meta"
$ws.Range("Q1").Value2 = "Flight_End_Date_meta"
$ws.Range("R1").Value2 = "Total_Budget_meta"
$ws.Range("S1").Value2 = "DSP_meta"
$ws.Range("T1").Value2 = "ML_Prediction"
# U1 (Budget_At_Risk) / V1 (ML_Early_Warning) unchanged.

# Drop the now-obsolete last column (old W: Excel_vs_ML_Disagree). This also
# removes its header cell and shrinks the sheet's dimension to A:V.
$ws.Columns.Item(23).Delete()

# ---------------------------------------------------------------------
# Sheet 3: Exec_Summary
# ---------------------------------------------------------------------
$es = $wb.Worksheets.Item("Exec_Summary")

# Remove the "Excel vs ML Disagreement Count" row; LAST_REFRESH_UTC shifts
# up from row 5 to row 4.
$es.Rows.Item(4).Delete()

# Refresh the timestamp in place.
$es.Range("B4").Value2 = "2026-02-16 06:09 UTC"
